$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new value in C2
$ws.Range("C2").Value = 10

# Update active selection to C2
$ws.Range("C2").Select()
